# Add data for 2021-11-17
# Updates the "Carjacking arrests by month, year over year" workbook:
#   - renames the sheet / title from "Through 2021-11-08" to "Through 2021-11-09"
#   - updates row 7 (May, 2021 columns) with a corrected count
#   - updates row 13 (November, partial month through 11-09) with new counts/rates
#   - updates row 14 (Total) with the new aggregated counts/rates

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet tab to match the new "through" date
$ws.Name = "Through 2021-11-09"

# Row 7 - May, 2021 arrest/no-arrest/rate shifted by one case
$ws.Range("T7").Value = 15
$ws.Range("U7").Value = 93
$ws.Range("V7").Value = 0.1389

# Row 13 - November (through 11-09), label + updated figures for every year column
$ws.Range("A13").Value = "November (through 11-09)"

$ws.Range("C13").Value = 12
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 21
$ws.Range("G13").Value = 0.0455
$ws.Range("I13").Value = 31
$ws.Range("J13").Value = 0.0312
$ws.Range("K13").Value = 5
$ws.Range("L13").Value = 18
$ws.Range("M13").Value = 0.2174
$ws.Range("N13").Value = 2
$ws.Range("P13").Value = 0.1538
$ws.Range("R13").Value = 56
$ws.Range("S13").Value = 0.0175
$ws.Range("U13").Value = 63
$ws.Range("V13").Value = 0.0156

# Row 14 - Total, updated figures for every year column
$ws.Range("C14").Value = 238
$ws.Range("D14").Value = 0.1185
$ws.Range("E14").Value = 53
$ws.Range("F14").Value = 455
$ws.Range("G14").Value = 0.1043
$ws.Range("I14").Value = 680
$ws.Range("J14").Value = 0.0836
$ws.Range("K14").Value = 71
$ws.Range("L14").Value = 567
$ws.Range("M14").Value = 0.1113
$ws.Range("N14").Value = 50
$ws.Range("P14").Value = 0.101
$ws.Range("R14").Value = 1059
$ws.Range("S14").Value = 0.0494
$ws.Range("T14").Value = 86
$ws.Range("U14").Value = 1422
$ws.Range("V14").Value = 0.057
